$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a string value while preventing Excel from
# auto-converting numeric-looking text (e.g. "506.73") into a float,
# which would corrupt trailing zeros / add floating point noise.
# We flip to text format ("@") for the assignment, then ClearFormats()
# so the cell keeps its original (default) style afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Update price (D) and volume-1h (E) columns for rows with changed values
$ws.Range("D2").Value = "57.719.97"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.439.57"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "506.73"
$ws.Range("E5").Value = "  -1.96%  "
Set-TextValue $ws.Range("D6") "129.31"
$ws.Range("E6").Value = "  -1.73%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "2.452.45"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "2.874.60"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "57.697.05"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "2.448.64"
$ws.Range("E18").Value = "  -0.67%  "
Set-TextValue $ws.Range("D19") "10.48"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  -0.95%  "
Set-TextValue $ws.Range("D21") "315.37"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  +0.04%  "
Set-TextValue $ws.Range("D23") "5.67"
$ws.Range("E23").Value = "  -1.14%  "
Set-TextValue $ws.Range("D24") "63.40"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  -0.11%  "
Set-TextValue $ws.Range("D26") "0.994"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("E28").Value = "  -1.24%  "
Set-TextValue $ws.Range("D29") "169.95"
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -2.73%  "
Set-TextValue $ws.Range("D31") "6.26"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.12%  "
Set-TextValue $ws.Range("D36") "17.73"
$ws.Range("E36").Value = "  -1.89%  "
Set-TextValue $ws.Range("D37") "1.26"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("E38").Value = "  -0.06%  "
Set-TextValue $ws.Range("D39") "36.31"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -1.56%  "
Set-TextValue $ws.Range("D41") "0.768"
$ws.Range("E41").Value = "  -2.36%  "
Set-TextValue $ws.Range("D42") "272.40"
$ws.Range("E42").Value = "  -0.41%  "

# Rows 43 and 44 swap places: Filecoin <-> RenderToken (rank order changed),
# each with freshly updated price and volume figures
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "5.01"
$ws.Range("E43").Value = "  +1.76%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "3.39"
$ws.Range("E44").Value = "  -2.29%  "

Set-TextValue $ws.Range("D45") "0.580"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  +0.18%  "
Set-TextValue $ws.Range("D47") "120.22"
$ws.Range("E47").Value = "  -5.07%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  -2.58%  "
Set-TextValue $ws.Range("D50") "0.0209"
$ws.Range("E50").Value = "  -2.08%  "
Set-TextValue $ws.Range("D51") "16.68"
$ws.Range("E51").Value = "  -1.83%  "
